# Update data.xlsx from the QR tool:
# Insert a new record as the second row (right after the header row),
# pushing all existing data rows down by one.
#
# Columns: id, code, name, address, mapsUrl, createdAt, note, phone,
#          branch, cccd, customerCode, officer, pinSalt, pinHash
# The new scan has no note/phone/branch/cccd/customerCode yet, so those
# cells are simply left blank (same as how Excel itself represents an
# empty cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("locations")

# Insert a new blank row at row 2 (shifts existing rows 2.. down to 3..)
$ws.Rows.Item(2).Insert()

# Fill in the new row with the freshly scanned QR record.
$ws.Cells.Item(2, 1).Value = "3ojnubbo1h6"
$ws.Cells.Item(2, 2).Value = "jh64ga9d"
$ws.Cells.Item(2, 3).Value = "a"
$ws.Cells.Item(2, 4).Value = "Phường Tăng Nhơn Phú, Thành phố Hồ Chí Minh, 71300, Việt Nam"
$ws.Cells.Item(2, 5).Value = "https://www.google.com/maps/search/?api=1&query=10.839061,106.792777"
$ws.Cells.Item(2, 6).Value = "2025-08-22T09:34:12.618Z"
$ws.Cells.Item(2, 12).Value = "Phan Minh Khải"
$ws.Cells.Item(2, 13).Value = "c5e29bf64454a22f"
$ws.Cells.Item(2, 14).Value = "2b93044ba18dd5aab233797be33ff611a3ddd62f00c9bd241013c86ffff2ae4c"
